# Update "想去人数" (F column) values across sheets per upstream data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 621
$ws1.Range("F5").Value = 4551
$ws1.Range("F6").Value = 1824
$ws1.Range("F7").Value = 125
$ws1.Range("F9").Value = 3074
$ws1.Range("F12").Value = 243
$ws1.Range("F13").Value = 593
$ws1.Range("F14").Value = 506
$ws1.Range("F15").Value = 509
$ws1.Range("F16").Value = 356
$ws1.Range("F17").Value = 129
$ws1.Range("F18").Value = 1758
$ws1.Range("F19").Value = 1302
$ws1.Range("F20").Value = 115
$ws1.Range("F21").Value = 1545
$ws1.Range("F22").Value = 125
$ws1.Range("F24").Value = 42
$ws1.Range("F25").Value = 524
$ws1.Range("F27").Value = 42
$ws1.Range("F31").Value = 3512
$ws1.Range("F32").Value = 742
$ws1.Range("F33").Value = 63
$ws1.Range("F34").Value = 243
$ws1.Range("F36").Value = 1697

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 37

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 621
$ws4.Range("F5").Value = 4551
$ws4.Range("F6").Value = 1824
$ws4.Range("F7").Value = 125
$ws4.Range("F9").Value = 3074
$ws4.Range("F12").Value = 243
$ws4.Range("F13").Value = 593
$ws4.Range("F14").Value = 506
$ws4.Range("F15").Value = 509
$ws4.Range("F17").Value = 356
$ws4.Range("F18").Value = 129
$ws4.Range("F19").Value = 1758
$ws4.Range("F20").Value = 1302
$ws4.Range("F21").Value = 115
$ws4.Range("F22").Value = 1545
$ws4.Range("F23").Value = 125
$ws4.Range("F25").Value = 42
$ws4.Range("F26").Value = 524
$ws4.Range("F28").Value = 42
$ws4.Range("F32").Value = 3514
$ws4.Range("F33").Value = 37
$ws4.Range("F34").Value = 742
$ws4.Range("F35").Value = 63
$ws4.Range("F36").Value = 243
$ws4.Range("F38").Value = 1697
